$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 1020
$wsExhibit.Range("F6").Value = 2384

# Sheet "全部类型" (rId4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 1020
$wsAll.Range("F8").Value = 2384
